$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.009099999999997
$ws.Range("D4").Value = -6.795599999999996
$ws.Range("C11").Value = -13.39029999999999
$ws.Range("C12").Value = -11.01849999999999
$ws.Range("D14").Value = -7.251000000000005
$ws.Range("C15").Value = -13.3561
$ws.Range("D26").Value = -8.171999999999999
$ws.Range("C27").Value = -13.3218
$ws.Range("C28").Value = -13.3852
$ws.Range("C31").Value = -13.7482
$ws.Range("D31").Value = -8.847400000000006
$ws.Range("C32").Value = -13.86910000000001
$ws.Range("D35").Value = -8.583099999999991
$ws.Range("C36").Value = -13.11650000000001
$ws.Range("D37").Value = -7.636799999999999
$ws.Range("C38").Value = -12.8472
$ws.Range("D39").Value = -7.996099999999999
$ws.Range("D40").Value = -8.026599999999991
$ws.Range("D45").Value = -7.551699999999999
$ws.Range("C46").Value = -14.68139999999999
$ws.Range("D52").Value = -7.508299999999996
$ws.Range("C54").Value = -13.11630000000001
$ws.Range("C55").Value = -13.75310000000001
$ws.Range("C56").Value = -13.38069999999999
$ws.Range("D57").Value = -8.414700000000003
$ws.Range("C67").Value = -10.74370000000001
$ws.Range("C69").Value = -11.6405
$ws.Range("C72").Value = -11.72930000000001
$ws.Range("C73").Value = -12.38260000000001
$ws.Range("D81").Value = -6.897299999999994
$ws.Range("C83").Value = -14.00820000000001
$ws.Range("D83").Value = -8.861099999999997
$ws.Range("C86").Value = -13.99049999999999
$ws.Range("C91").Value = -10.3557
$ws.Range("C93").Value = -10.9433
$ws.Range("C99").Value = -13.2845
$ws.Range("D100").Value = -8.079000000000004
$ws.Range("D102").Value = -7.638
